$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new poll row (row 4): "Test Question 3" / "Textbox" type, deactivated,
# with no option columns (C:G) since a Textbox question has no options.

# Copy row 3's formatting down to row 4 first, so the new row picks up the
# same look (font/number-format) as the existing data rows.
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I4").PasteSpecial()
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("A4").Value = "Test Question 3"
$ws.Range("B4").Value = "Textbox"
$ws.Range("C4:G4").ClearContents()
$ws.Range("H4").Value = "deactivate"

# I4 should hold the literal text "true" (matching I2), not a boolean, so
# copy it over from I2 to preserve the text formatting/type.
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("H4").Select() | Out-Null
